# Weekly update: insert a new record as the new first row of this week's
# batch (row 51), pushing the existing rows 51-76 down to 52-77.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51; existing row 51 (and below) shift down to 52.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new Espinaca price record.
$ws.Range('A51').Value = 10
$ws.Range('B51').Value = 'Vega Modelo de Temuco'
$ws.Range('C51').Value = 'La Araucanía'
$ws.Range('D51').Value = 44455
$ws.Range('E51').Value = 9
$ws.Range('F51').Value = 100112012
$ws.Range('G51').Value = 'Espinaca'
$ws.Range('H51').Value = 'Sin especificar'
$ws.Range('I51').Value = 'Primera'
$ws.Range('J51').Value = 10
$ws.Range('K51').Value = 9000
$ws.Range('L51').Value = 9000
$ws.Range('M51').Value = 9000
$ws.Range('N51').Value = '$/docena de atados'
$ws.Range('O51').Value = 'Región de La Araucanía'
$ws.Range('P51').Value = 3000
$ws.Range('Q51').Value = 3
$ws.Range('R51').Value = 'Hortaliza'
